$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(
    "2024-10-05",
    "2024-10-04",
    "2024-10-03",
    "2024-10-02",
    "2024-10-01",
    "2024-09-29",
    "2024-09-27",
    "2024-09-23",
    "2024-09-21",
    "2024-09-19",
    "2024-09-18",
    "2024-09-16",
    "2024-09-14",
    "2024-09-13",
    "2024-09-12",
    "2024-09-11",
    "2024-09-10",
    "2024-09-09",
    "2024-09-08",
    "2024-09-07",
    "2023-11-10",
    "2024-01-09"
)

$startRow = 37
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
}
